$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = -7.816199999999993
$ws.Range("D12").Value = -5.915999999999997
$ws.Range("D15").Value = -8.373600000000003
$ws.Range("D27").Value = -7.876700000000004
$ws.Range("D28").Value = -7.832100000000001
$ws.Range("D31").Value = -7.403899999999997
$ws.Range("D32").Value = -7.426099999999992
$ws.Range("D36").Value = -7.871400000000005
$ws.Range("D38").Value = -7.33
$ws.Range("D46").Value = -7.874799999999996
$ws.Range("D54").Value = -8.053200000000006
$ws.Range("D55").Value = -7.393599999999993
$ws.Range("D56").Value = -8.957300000000005
$ws.Range("D67").Value = -7.313099999999999
$ws.Range("D69").Value = -7.209599999999998
$ws.Range("D72").Value = -7.046999999999998
$ws.Range("D73").Value = -7.865399999999994
$ws.Range("D83").Value = -8.434000000000008
$ws.Range("D86").Value = -8.5052
$ws.Range("D91").Value = -7.849200000000003
$ws.Range("D93").Value = -6.993299999999993
$ws.Range("D99").Value = -7.928899999999999
$ws.Range("D104").Value = -7.552100000000004
$ws.Range("D105").Value = -8.135800000000003
